# Auto-generated edit script: update market-price derived cells per scheduled runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 646.44446
$ws.Range("J33").Value = 387.75
$ws.Range("L33").Value = 387.75
$ws.Range("N33").Value = -845.75

$ws.Range("H43").Value = 4155
$ws.Range("J43").Value = 5500.143
$ws.Range("L43").Value = 5500.143
$ws.Range("N43").Value = -5638.143

$ws.Range("H86").Value = 9804.904
$ws.Range("I86").Value = 7000.8
$ws.Range("J86").Value = 12354.091
$ws.Range("K86").Value = 7000.8
$ws.Range("L86").Value = 12354.091
$ws.Range("M86").Value = -5877.8
$ws.Range("N86").Value = -14600.091

$ws.Range("H89").Value = 9804.904
$ws.Range("I89").Value = 7000.8
$ws.Range("J89").Value = 12354.091
$ws.Range("K89").Value = 35004
$ws.Range("L89").Value = 61770.455
$ws.Range("M89").Value = -29388
$ws.Range("N89").Value = -73002.455

$ws.Range("H132").Value = 26318470
$ws.Range("I132").Value = 26318470
$ws.Range("K132").Value = 78955410
$ws.Range("M132").Value = -78952880

$ws.Range("H137").Value = 76838.086
$ws.Range("I137").Value = 113388.5
$ws.Range("J137").Value = 3737.25
$ws.Range("K137").Value = 340165.5
$ws.Range("L137").Value = 11211.75
$ws.Range("M137").Value = -337615.5
$ws.Range("N137").Value = -16311.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2420.21
$ws.Range("I32").Value = 1671.1428
$ws.Range("J32").Value = 9994.111000000001
$ws.Range("K32").Value = 1671.1428
$ws.Range("L32").Value = 9994.111000000001
$ws.Range("M32").Value = -1384.1428
$ws.Range("N32").Value = -10568.111

$ws.Range("H61").Value = 2906.5881
$ws.Range("I61").Value = 2494.0715
$ws.Range("K61").Value = 2494.0715
$ws.Range("M61").Value = -2282.0715

$ws.Range("H136").Value = 2906.5881
$ws.Range("I136").Value = 2494.0715
$ws.Range("K136").Value = 7482.2145
$ws.Range("M136").Value = -4932.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6265818.5
$ws.Range("I86").Value = 7702623
$ws.Range("K86").Value = 7702623
$ws.Range("M86").Value = -7701500

$ws.Range("H89").Value = 6265818.5
$ws.Range("I89").Value = 7702623
$ws.Range("K89").Value = 38513115
$ws.Range("M89").Value = -38507499

$ws.Range("H99").Value = 6804909.5
$ws.Range("I99").Value = 10205580
$ws.Range("J99").Value = 3569.1428
$ws.Range("K99").Value = 10205580
$ws.Range("L99").Value = 3569.1428
$ws.Range("M99").Value = -10204082
$ws.Range("N99").Value = -6565.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7998.8887
$ws.Range("I58").Value = 10631.272
$ws.Range("J58").Value = 3862.2856
$ws.Range("K58").Value = 10631.272
$ws.Range("L58").Value = 3862.2856
$ws.Range("M58").Value = -10428.272
$ws.Range("N58").Value = -4268.2856

$ws.Range("H81").Value = 93280
$ws.Range("J81").Value = 93280
$ws.Range("L81").Value = 93280
$ws.Range("N81").Value = -95276

$ws.Range("H84").Value = 93280
$ws.Range("J84").Value = 93280
$ws.Range("L84").Value = 279840
$ws.Range("N84").Value = -289824

$ws.Range("H99").Value = 5873.75
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 5998.3335
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 5998.3335
$ws.Range("M99").Value = -4002
$ws.Range("N99").Value = -8994.333500000001

$ws.Range("H126").Value = 5873.75
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 5998.3335
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 17995.0005
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -22935.0005

$ws.Range("H132").Value = 72476.25999999999
$ws.Range("I132").Value = 79377.62
$ws.Range("J132").Value = 66067.86
$ws.Range("K132").Value = 238132.86
$ws.Range("L132").Value = 198203.58
$ws.Range("M132").Value = -235602.86
$ws.Range("N132").Value = -203263.58

$ws.Range("H134").Value = 35584.145
$ws.Range("I134").Value = 49671.05
$ws.Range("K134").Value = 149013.15
$ws.Range("M134").Value = -146478.15

$ws.Range("H136").Value = 7998.8887
$ws.Range("I136").Value = 10631.272
$ws.Range("J136").Value = 3862.2856
$ws.Range("K136").Value = 31893.816
$ws.Range("L136").Value = 11586.8568
$ws.Range("M136").Value = -29343.816
$ws.Range("N136").Value = -16686.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 210717.31
$ws.Range("J46").Value = 3390.3635
$ws.Range("L46").Value = 10171.0905
$ws.Range("N46").Value = -10353.0905

$ws.Range("H50").Value = 1096.909
$ws.Range("I50").Value = 845.8333
$ws.Range("J50").Value = 1398.2
$ws.Range("K50").Value = 2537.4999
$ws.Range("L50").Value = 4194.6
$ws.Range("M50").Value = -2056.4999
$ws.Range("N50").Value = -5156.6

$ws.Range("H53").Value = 1096.909
$ws.Range("I53").Value = 845.8333
$ws.Range("J53").Value = 1398.2
$ws.Range("K53").Value = 2537.4999
$ws.Range("L53").Value = 4194.6
$ws.Range("M53").Value = -2056.4999
$ws.Range("N53").Value = -5156.6

$ws.Range("H122").Value = 1033.6428
$ws.Range("J122").Value = 846
$ws.Range("L122").Value = 7614
$ws.Range("N122").Value = -12514

$ws.Range("H129").Value = 1041.7778
$ws.Range("I129").Value = 1041.7778
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3125.3334
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1874.6666
$ws.Range("N129").ClearContents()

$ws.Range("H131").Value = 13444520
$ws.Range("I131").Value = 8335848.5
$ws.Range("J131").Value = 15877221
$ws.Range("K131").Value = 25007545.5
$ws.Range("L131").Value = 47631663
$ws.Range("M131").Value = -25002505.5
$ws.Range("N131").Value = -47641743

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2497.8936
$ws.Range("I132").Value = 2386.9333
$ws.Range("J132").Value = 4994.5
$ws.Range("K132").Value = 7160.7999
$ws.Range("L132").Value = 14983.5
$ws.Range("M132").Value = -4630.7999
$ws.Range("N132").Value = -20043.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4726.923
$ws.Range("I40").Value = 2859.2727
$ws.Range("J40").Value = 14999
$ws.Range("K40").Value = 2859.2727
$ws.Range("L40").Value = 14999
$ws.Range("M40").Value = -2723.2727
$ws.Range("N40").Value = -15271

$ws.Range("H46").Value = 4768.8096
$ws.Range("J46").Value = 7528.75
$ws.Range("L46").Value = 7528.75
$ws.Range("N46").Value = -7904.75

$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 500
$ws.Range("K68").Value = 500
$ws.Range("M68").Value = 249

$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 500
$ws.Range("K71").Value = 2500
$ws.Range("M71").Value = 1244

$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -22350

$ws.Range("H122").Value = 5127.476
$ws.Range("I122").Value = 3371
$ws.Range("J122").Value = 6005.7144
$ws.Range("K122").Value = 10113
$ws.Range("L122").Value = 18017.1432
$ws.Range("M122").Value = -7663
$ws.Range("N122").Value = -22917.1432

$ws.Range("H132").Value = 7622.7812
$ws.Range("I132").Value = 7404.6924
$ws.Range("K132").Value = 22214.0772
$ws.Range("M132").Value = -19684.0772

$ws.Range("H136").Value = 52555.88
$ws.Range("I136").Value = 82307.60000000001
$ws.Range("J136").Value = 6068.8125
$ws.Range("K136").Value = 246922.8
$ws.Range("L136").Value = 18206.4375
$ws.Range("M136").Value = -244372.8
$ws.Range("N136").Value = -23306.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 15115.223
$ws.Range("I31").Value = 9999.666999999999
$ws.Range("J31").Value = 17673
$ws.Range("K31").Value = 9999.666999999999
$ws.Range("L31").Value = 17673
$ws.Range("M31").Value = -9651.666999999999
$ws.Range("N31").Value = -18369

$ws.Range("H122").Value = 2882.276
$ws.Range("I122").Value = 2759.7778
$ws.Range("J122").Value = 3082.7273
$ws.Range("K122").Value = 8279.3334
$ws.Range("L122").Value = 9248.1819
$ws.Range("M122").Value = -5829.3334
$ws.Range("N122").Value = -14148.1819
